# khl/Injuries_Master_Clubs.xlsx — publish files + archive (2025-12-01 15:05:12)
#
# 1) "snapshot" sheet: refresh the scraped_at (column K) timestamps for rows 2-30
#    with the latest scrape run's values.
# 2) "new_injured" sheet: the single staged row (row 2, Трактор / Мыльников Сергей И)
#    was consumed/published, so remove it — dimension shrinks back to A1:G1.

$wb = $excel.ActiveWorkbook

# --- 1) Update scraped_at (column K) on the "snapshot" sheet ---------------
$snapshot = $wb.Worksheets.Item("snapshot")

$snapshot.Cells.Item(2, 11).Value  = "2025-12-01T07:01:43.019771+00:00"
$snapshot.Cells.Item(3, 11).Value  = "2025-12-01T07:01:45.315650+00:00"
$snapshot.Cells.Item(4, 11).Value  = "2025-12-01T07:01:45.315681+00:00"
$snapshot.Cells.Item(5, 11).Value  = "2025-12-01T07:01:48.063419+00:00"
$snapshot.Cells.Item(6, 11).Value  = "2025-12-01T07:01:50.835017+00:00"
$snapshot.Cells.Item(7, 11).Value  = "2025-12-01T07:01:53.655587+00:00"
$snapshot.Cells.Item(8, 11).Value  = "2025-12-01T07:01:55.872696+00:00"
$snapshot.Cells.Item(9, 11).Value  = "2025-12-01T07:02:01.133972+00:00"
$snapshot.Cells.Item(10, 11).Value = "2025-12-01T07:02:01.134006+00:00"
$snapshot.Cells.Item(11, 11).Value = "2025-12-01T07:02:03.902337+00:00"
$snapshot.Cells.Item(12, 11).Value = "2025-12-01T07:02:06.185439+00:00"
$snapshot.Cells.Item(13, 11).Value = "2025-12-01T07:02:06.185470+00:00"
$snapshot.Cells.Item(14, 11).Value = "2025-12-01T07:02:09.679960+00:00"
$snapshot.Cells.Item(15, 11).Value = "2025-12-01T07:02:12.012664+00:00"
$snapshot.Cells.Item(16, 11).Value = "2025-12-01T07:02:12.012691+00:00"
$snapshot.Cells.Item(17, 11).Value = "2025-12-01T07:02:12.012709+00:00"
$snapshot.Cells.Item(18, 11).Value = "2025-12-01T07:02:14.338200+00:00"
$snapshot.Cells.Item(19, 11).Value = "2025-12-01T07:02:14.338233+00:00"
$snapshot.Cells.Item(20, 11).Value = "2025-12-01T07:02:16.605064+00:00"
$snapshot.Cells.Item(21, 11).Value = "2025-12-01T07:02:16.605094+00:00"
$snapshot.Cells.Item(22, 11).Value = "2025-12-01T07:02:18.904705+00:00"
$snapshot.Cells.Item(23, 11).Value = "2025-12-01T07:02:18.904743+00:00"
$snapshot.Cells.Item(24, 11).Value = "2025-12-01T07:02:18.904765+00:00"
$snapshot.Cells.Item(25, 11).Value = "2025-12-01T07:02:21.211591+00:00"
$snapshot.Cells.Item(26, 11).Value = "2025-12-01T07:02:29.844045+00:00"
$snapshot.Cells.Item(27, 11).Value = "2025-12-01T07:02:32.104455+00:00"
$snapshot.Cells.Item(28, 11).Value = "2025-12-01T07:02:32.104485+00:00"
$snapshot.Cells.Item(29, 11).Value = "2025-12-01T07:02:34.443683+00:00"
$snapshot.Cells.Item(30, 11).Value = "2025-12-01T07:02:34.443713+00:00"

# --- 2) Remove the now-stale staged row on "new_injured" -------------------
$newInjured = $wb.Worksheets.Item("new_injured")
$newInjured.Rows.Item(2).Delete()
